$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2025-07-12"

$arr = New-Object "object[,]" 51,4
$arr[0,0] = 'rank'
$arr[0,1] = 'title'
$arr[0,2] = 'author'
$arr[0,3] = 'latest_episode'

$arr[1,0] = 1
$arr[1,1] = '魔術師クノンは見えている'
$arr[1,2] = 'La-na(作画) 南野海風(原作) Ｌａｒｕｈａ(キャラクター原案)'
$arr[1,3] = '第38話②'
$arr[2,0] = 2
$arr[2,1] = '悪人面したＢ級冒険者 主人公とその幼馴染たちのパパになる'
$arr[2,2] = 'こげめ(著者) えんじ(原作) ハラカズヒロ(キャラクター原案)'
$arr[2,3] = '「名もなき英雄譚」前半'
$arr[3,0] = 3
$arr[3,1] = '勇者に全部奪われた俺は勇者の母親とパーティを組みました！'
$arr[3,2] = '久遠まこと(著者) 石のやっさん(原作)'
$arr[3,3] = '第28話'
$arr[4,0] = 4
$arr[4,1] = 'ダンジョンの幼なじみ'
$arr[4,2] = '久真やすひさ(著者)'
$arr[4,3] = '第55話'
$arr[5,0] = 5
$arr[5,1] = '金属スライムを倒しまくった俺が【黒鋼の王】と呼ばれるまで'
$arr[5,2] = '藤屋いずこ(著者) 温泉カピバラ(原作) 山椒魚(キャラクター原案)'
$arr[5,3] = '第13章-2'
$arr[6,0] = 6
$arr[6,1] = '淫獄団地'
$arr[6,2] = '搾精研究所(原作) 丈山雄為(漫画)'
$arr[6,3] = '第48話（後編）'
$arr[7,0] = 7
$arr[7,1] = 'まんきつしたい常連さん'
$arr[7,2] = 'しんみりん(著者)'
$arr[7,3] = '第45話後編'
$arr[8,0] = 8
$arr[8,1] = '女友達は頼めば意外とヤらせてくれる'
$arr[8,2] = 'ろくろ(漫画) 鏡遊(原作)'
$arr[8,3] = '特別イラスト'
$arr[9,0] = 9
$arr[9,1] = '美人女上司滝沢さん'
$arr[9,2] = 'やんBARU(著者)'
$arr[9,3] = '第201話'
$arr[10,0] = 10
$arr[10,1] = 'よくわからないけれど異世界に転生していたようです'
$arr[10,2] = '内々けやき あし カオミン'
$arr[10,3] = '第135話 よくわからないけれど導かれてしまったようです（２）'
$arr[11,0] = 11
$arr[11,1] = '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~'
$arr[11,2] = '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)'
$arr[11,3] = '第80話その2'
$arr[12,0] = 12
$arr[12,1] = '解雇された暗黒兵士(30代)のスローなセカンドライフ'
$arr[12,2] = '岡沢六十四 るれくちぇ sage・ジョー'
$arr[12,3] = '第70話(後編) 特別報酬ミッション'
$arr[13,0] = 13
$arr[13,1] = '老後に備えて異世界で８万枚の金貨を貯めます'
$arr[13,2] = 'FUNA 東西 モトエ恵介'
$arr[13,3] = '第119話　会談［その4］'
$arr[14,0] = 14
$arr[14,1] = '陰キャの俺が席替えでS級美少女に囲まれたら秘密の関係が始まった。'
$arr[14,2] = '星野 星野(原作) バラマツヒトミ(漫画) 黒兎 ゆう(キャラクター原案)'
$arr[14,3] = '第3話'
$arr[15,0] = 15
$arr[15,1] = 'ノロマ魔法と呼ばれた魔法使いは重力魔法で無双する　～まだ重力の概念のない世界にて、少年は万有引力の王となる～'
$arr[15,2] = '神原絵理華(漫画) 一森一輝(原作)'
$arr[15,3] = '第17話④'
$arr[16,0] = 16
$arr[16,1] = '不徳のギルド'
$arr[16,2] = '河添太一'
$arr[16,3] = '第９６話：分福'
$arr[17,0] = 17
$arr[17,1] = 'ホームセンターごと呼び出された私の大迷宮リノベーション！'
$arr[17,2] = 'ばたっち(漫画) 星崎崑(原作) 志田(キャラクター原案)'
$arr[17,3] = '第4話後編'
$arr[18,0] = 18
$arr[18,1] = '仕事帰り、独身の美人上司に頼まれて'
$arr[18,2] = '望公太(原作) とんのすけ(作画) しの(キャラクター原案)'
$arr[18,3] = '第19話-2'
$arr[19,0] = 19
$arr[19,1] = 'ダウナーお姉さんは遊びたい'
$arr[19,2] = '山鷹景'
$arr[19,3] = '第13話'
$arr[20,0] = 20
$arr[20,1] = '農学博士の異世界無双～禁忌の知識で築くモンスター娘ハーレム～'
$arr[20,2] = 'インド僧(原作) ヤスウミ(作画)'
$arr[20,3] = '第24話'
$arr[21,0] = 21
$arr[21,1] = '最強で最速の無限レベルアップ ～スキル【経験値1000倍】と【レベルフリー】でレベル上限の枷が外れた俺は無双する～'
$arr[21,2] = 'シオヤマ琴 鳥羽田 航 トモゼロ'
$arr[21,3] = '休載マンガ'
$arr[22,0] = 22
$arr[22,1] = 'センパイ、自宅警備員の雇用はいかがですか？'
$arr[22,2] = '漫画：コブラサナギ 原作：二上圭 キャラ原案：日向あずり'
$arr[22,3] = '第5話前半'
$arr[23,0] = 23
$arr[23,1] = '修羅幼女の英雄譚～半端者と言われた傭兵、幼女に転生して成り上がる～'
$arr[23,2] = '作画：むらたん 原作：沙城流'
$arr[23,3] = '第7話(1)'
$arr[24,0] = 24
$arr[24,1] = 'ある日突然、ギャルの許嫁ができた'
$arr[24,2] = '窪茶(漫画) 泉谷一樹(原作) なかむら(文庫イラスト) まめぇ(原作イラスト)'
$arr[24,3] = '第12話'
$arr[25,0] = 25
$arr[25,1] = '僕のいけずな婚約者'
$arr[25,2] = '冬谷リク(漫画)'
$arr[25,3] = '第7話'
$arr[26,0] = 26
$arr[26,1] = '異世界はスマートフォンとともに。'
$arr[26,2] = 'そと(漫画) 冬原パトラ(原作) 兎塚エイジ(キャラクター原案)'
$arr[26,3] = 'EPISODE:102‐②'
$arr[27,0] = 27
$arr[27,1] = 'スキルがなければレベルを上げる～９９がカンストの世界でレベル800万からスタート～'
$arr[27,2] = '倉橋ユウス(漫画) 岡沢六十四(原作)'
$arr[27,3] = '第51話②'
$arr[28,0] = 28
$arr[28,1] = 'ヤンデレかと思ったらもっとヤベー女だった'
$arr[28,2] = '八木戸マト'
$arr[28,3] = '第66話　最後に彼氏の全てが欲しいヤンデレ彼女'
$arr[29,0] = 29
$arr[29,1] = '時森さんが無防備です!!'
$arr[29,2] = 'たざわ'
$arr[29,3] = '第62話'
$arr[30,0] = 30
$arr[30,1] = '不老不死少女の苗床旅行記'
$arr[30,2] = 'ふじはん(漫画) ルナ・ウサギ(原作)'
$arr[30,3] = '第16話前編'
$arr[31,0] = 31
$arr[31,1] = 'ゴミ以下だと追放された使用人、実は前世賢者です　～史上最強の賢者、世界最高峰の学園に通う～'
$arr[31,2] = '原作：夜分長文 漫画：矢部利恩 キャラクター原案：蔓木鋼音'
$arr[31,3] = '第14話 魔女対策（１）'
$arr[32,0] = 32
$arr[32,1] = '辺境モブ貴族のウチに嫁いできた悪役令嬢が、めちゃくちゃできる良い嫁なんだが？'
$arr[32,2] = 'tera(原作) 朝倉はやて(作画) 徹田(キャラクター原案)'
$arr[32,3] = '第9話'
$arr[33,0] = 33
$arr[33,1] = '断れない会長は友江くんにだけしてあげたい'
$arr[33,2] = '沼地どろまる(著者)'
$arr[33,3] = '休載漫画'
$arr[34,0] = 34
$arr[34,1] = 'くじ引き特賞：無双ハーレム権'
$arr[34,2] = '原作／三木なずな（GA文庫／SBクリエイティブ刊） 漫画／長谷見亮 キャラクター原案／瑠奈璃亜'
$arr[34,3] = '第58話-01　新たな王女たちへ、受け継がれし慈愛の心！'
$arr[35,0] = 35
$arr[35,1] = '脱稿するまでオチません'
$arr[35,2] = 'ヨシラギ(著者)'
$arr[35,3] = '第32話前半'
$arr[36,0] = 36
$arr[36,1] = '義妹生活'
$arr[36,2] = '三河ごーすと(原作) 奏ユミカ(漫画) Hiten(キャラクター原案)'
$arr[36,3] = '第30話-1'
$arr[37,0] = 37
$arr[37,1] = 'クラスメイトは異世界で勇者になったけど、俺だけ現代日本に置き去りにされました'
$arr[37,2] = 'カボチャマスク(原作) 仲紙(漫画)'
$arr[37,3] = '第9話-4'
$arr[38,0] = 38
$arr[38,1] = '俺の愛娘は悪役令嬢'
$arr[38,2] = 'かわもり かぐら(原作) ほづみりや(漫画) 縞(キャラクター原案)'
$arr[38,3] = '第4話-1'
$arr[39,0] = 39
$arr[39,1] = '俺の『全自動支援（フルオートバフ）』で仲間たちが世界最強 ～そこにいるだけ無自覚無双～'
$arr[39,2] = 'IプルT(著者) epina(原作) 片倉響(キャラクター原案)'
$arr[39,3] = '第３３話「砂浜の盗賊たち」'
$arr[40,0] = 40
$arr[40,1] = '勇者パーティから追い出された不遇職【罠士】、ユニークスキル【矢印】で最強になる'
$arr[40,2] = '作画：たつひこ 原作：白石 有希'
$arr[40,3] = '第7話(1)'
$arr[41,0] = 41
$arr[41,1] = 'その冒険者、取り扱い注意。 ～正体は無敵の下僕たちを統べる異世界最強の魔導王～'
$arr[41,2] = '満月シオン(作画) Sin Guilty(ツギクル)(原作) M.B(キャラクター原案)'
$arr[41,3] = '56章　はじまりの愚か者②　前編'
$arr[42,0] = 42
$arr[42,1] = 'モブ高生の俺でも冒険者になればリア充になれますか？'
$arr[42,2] = '原作：百均 漫画：さぎやまれん キャラクター原案：hai'
$arr[42,3] = '第29.5話'
$arr[43,0] = 43
$arr[43,1] = 'クロの戦記Ⅱ 異世界転移した僕が最強なのはベッドの上だけのようです'
$arr[43,2] = 'サイトウアユム(原作) ユリシロ(漫画) むつみまさと(キャラクター原案)'
$arr[43,3] = '第22話-1'
$arr[44,0] = 44
$arr[44,1] = 'ぽんドロイド！ はまさん'
$arr[44,2] = 'はれやまはれぞう(著者)'
$arr[44,3] = '第3話'
$arr[45,0] = 45
$arr[45,1] = 'ギルドを追放された回復術士、実は魔力無限だったので規格外の回復魔法で伝説となる'
$arr[45,2] = '漫画：坂下コウ 原作：霞杏檎'
$arr[45,3] = '第3話(3)'
$arr[46,0] = 46
$arr[46,1] = '聖剣が最強の世界で、少年は弓に愛される～封印された魔王がくれた力で聖剣士たちを援護します～'
$arr[46,2] = 'さとう(原作) 貞清カズヒコ(漫画)'
$arr[46,3] = '第13話①'
$arr[47,0] = 47
$arr[47,1] = 'おじ転生'
$arr[47,2] = '相葉キョウコ'
$arr[47,3] = '第14話'
$arr[48,0] = 48
$arr[48,1] = '俺以外誰も採取できない素材なのに「素材採取率が低い」とパワハラする幼馴染錬金術師と絶縁した専属魔導士、辺境の町でスローライフを送りたい。'
$arr[48,2] = '狐御前(原作) 西岡知三(作画) ＮＯＣＯ(キャラクター原案)'
$arr[48,3] = '第23話-2'
$arr[49,0] = 49
$arr[49,1] = '落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～'
$arr[49,2] = '村上よしゆき 茨木野 あるてら'
$arr[49,3] = '第４０話　勇者、聖女と元聖騎士と再会し、魚人を追っ払う（４）'
$arr[50,0] = 50
$arr[50,1] = '最も嫌われている最凶の悪役に転生'
$arr[50,2] = '灰色の鼠(原作) 沢田かに(漫画)'
$arr[50,3] = '第17話①'

$fullRange = $newSheet.Range("A1:D51")
$fullRange.Value = $arr

$headerRange = $newSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

Write-Output "done"
